# "Update today hours. Again"
#
# The timesheet tracks planned vs. actual working hours for two weeks.
# Row 14/15 (merged) is the "Фактические часы работы" (actual hours
# worked) line for the Friday of the second week: K14 holds the
# clock-in time and K15 the clock-out time; L14 is the computed hours
# worked that day, and M14 is the week's running total.
#
# Today's clock-out time moved an hour earlier (19:45 -> 18:45), which
# shortens Friday's duration by an hour and the weekly total along
# with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clock-out time for Friday (actual hours row) - was 19:45
$ws.Range("K15").Value = 0.78125          # 18:45

# Recomputed hours worked that day (was 4.25)
$ws.Range("L14").Value = 3.25

# Recomputed weekly total (was 21)
$ws.Range("M14").Value = 20
